# New crime data collected - weekly CompStat figures updated (week of
# 10/23/2023 through 10/29/2023, Volume 30 Number 43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/number and the reporting week's date range -----------
# These are rich-text shared strings made of several runs that all share the
# same font, so a plain text replace reproduces the same visible content.
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Cells whose value flips between a number and the literal text "0" -----
# (blank-week precincts are shown as text "0" instead of a numeric 0; when a
# precinct goes from "no incidents" to "some incidents" or vice versa the
# cell's type - and therefore its style - has to change with it.)

function Set-TextZero($ref) {
    # Force a literal text entry (leading apostrophe => text, not number),
    # then copy the neighboring "text zero" cell's number format onto it so
    # the resulting style matches the rest of the "0" cells in the sheet.
    $ws.Range($ref).Value = "'0"
    $ws.Range("D14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-Number($ref, $value) {
    # Force a numeric entry and copy a plain numeric sibling's style onto it.
    $ws.Range($ref).Value = $value
    $ws.Range("D15").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

Set-TextZero "C14"
Set-Number   "C15" 1
Set-TextZero "C22"
Set-Number   "C26" 1
Set-TextZero "C28"
Set-TextZero "C29"

# --- Remaining weekly/28-day/YTD counts and their computed % changes -------
$ws.Range("L14").Value = -26.666666666666
$ws.Range("N14").Value = -81.666666666666
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = -20
$ws.Range("I15").Value = 38
$ws.Range("J15").Value = 41
$ws.Range("K15").Value = -7.317073170731
$ws.Range("L15").Value = 8.571428571428
$ws.Range("M15").Value = 22.58064516129
$ws.Range("N15").Value = -39.682539682539
$ws.Range("D16").Value = 19
$ws.Range("E16").Value = -52.631578947368
$ws.Range("G16").Value = 55
$ws.Range("H16").Value = -41.818181818181
$ws.Range("I16").Value = 446
$ws.Range("J16").Value = 615
$ws.Range("K16").Value = -27.479674796748
$ws.Range("L16").Value = 11.779448621553
$ws.Range("M16").Value = 5.687203791469
$ws.Range("N16").Value = -73.841642228739
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = -10.526315789473
$ws.Range("I17").Value = 708
$ws.Range("J17").Value = 660
$ws.Range("K17").Value = 7.272727272727
$ws.Range("L17").Value = 21.440823327615
$ws.Range("M17").Value = 64.651162790697
$ws.Range("N17").Value = -11.610486891385
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 262
$ws.Range("J18").Value = 262
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 23.584905660377
$ws.Range("M18").Value = -20.121951219512
$ws.Range("N18").Value = -82.931596091205
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 81
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 760
$ws.Range("J19").Value = 864
$ws.Range("K19").Value = -12.037037037037
$ws.Range("L19").Value = 23.778501628664
$ws.Range("M19").Value = 60.337552742616
$ws.Range("N19").Value = 17.283950617283
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 39
$ws.Range("H20").Value = -7.142857142857
$ws.Range("I20").Value = 544
$ws.Range("J20").Value = 442
$ws.Range("K20").Value = 23.076923076923
$ws.Range("L20").Value = 56.772334293948
$ws.Range("M20").Value = 157.81990521327
$ws.Range("N20").Value = -67.405632115038
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 68
$ws.Range("E21").Value = -30.882352941176
$ws.Range("F21").Value = 196
$ws.Range("G21").Value = 267
$ws.Range("H21").Value = -26.591760299625
$ws.Range("I21").Value = 2769
$ws.Range("J21").Value = 2890
$ws.Range("K21").Value = -4.186851211072
$ws.Range("L21").Value = 25.578231292517
$ws.Range("M21").Value = 45.354330708661
$ws.Range("N21").Value = -57.275111865452
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -27.777777777777
$ws.Range("M22").Value = -7.142857142857
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -42.857142857142
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = -8.333333333333
$ws.Range("I23").Value = 261
$ws.Range("J23").Value = 279
$ws.Range("K23").Value = -6.451612903225
$ws.Range("L23").Value = 17.567567567567
$ws.Range("M23").Value = 38.829787234042
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 3.571428571428
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 156
$ws.Range("H24").Value = -23.076923076923
$ws.Range("I24").Value = 1616
$ws.Range("J24").Value = 1706
$ws.Range("K24").Value = -5.2754982415
$ws.Range("L24").Value = 38.712446351931
$ws.Range("M24").Value = 29.487179487179
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 65
$ws.Range("G25").Value = 94
$ws.Range("H25").Value = -30.851063829787
$ws.Range("I25").Value = 934
$ws.Range("J25").Value = 940
$ws.Range("K25").Value = -0.63829787234
$ws.Range("L25").Value = 15.880893300248
$ws.Range("M25").Value = -29.0273556231
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 54
$ws.Range("J26").Value = 61
$ws.Range("K26").Value = -11.475409836065
$ws.Range("L26").Value = -6.896551724137
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 95
$ws.Range("J27").Value = 59
$ws.Range("K27").Value = 61.016949152542
$ws.Range("L27").Value = 9.19540229885
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("L28").Value = -28.301886792452
$ws.Range("M28").Value = -7.317073170731
$ws.Range("N28").Value = -70.542635658914
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 300
$ws.Range("L29").Value = -27.272727272727
$ws.Range("M29").Value = -8.571428571428
$ws.Range("N29").Value = -72.173913043478
